$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4238.2856
$ws.Range("I32").Value = 2399
$ws.Range("J32").Value = 4544.8335
$ws.Range("K32").Value = 2399
$ws.Range("L32").Value = 4544.8335
$ws.Range("M32").Value = -2073
$ws.Range("N32").Value = -5196.8335

$ws.Range("H98").Value = 5418.3794
$ws.Range("I98").Value = 4542.3
$ws.Range("J98").Value = 7365.222
$ws.Range("K98").Value = 4542.3
$ws.Range("L98").Value = 7365.222
$ws.Range("M98").Value = -3044.3
$ws.Range("N98").Value = -10361.222

$ws.Range("H122").Value = 5418.3794
$ws.Range("I122").Value = 4542.3
$ws.Range("J122").Value = 7365.222
$ws.Range("K122").Value = 13626.9
$ws.Range("L122").Value = 22095.666
$ws.Range("M122").Value = -11176.9
$ws.Range("N122").Value = -26995.666

$ws.Range("H137").Value = 1956
$ws.Range("I137").Value = 1971.25
$ws.Range("J137").Value = 1895
$ws.Range("K137").Value = 5913.75
$ws.Range("L137").Value = 5685
$ws.Range("M137").Value = -3363.75
$ws.Range("N137").Value = -10785

$ws.Range("H140").Value = 139999.5
$ws.Range("I140").Value = 130000
$ws.Range("J140").Value = 149999
$ws.Range("K140").Value = 130000
$ws.Range("L140").Value = 149999
$ws.Range("M140").Value = -124820
$ws.Range("N140").Value = -160359

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 31251492
$ws.Range("I74").Value = 35715260
$ws.Range("J74").Value = 5124.25
$ws.Range("K74").Value = 35715260
$ws.Range("L74").Value = 5124.25
$ws.Range("M74").Value = -35714386
$ws.Range("N74").Value = -6872.25

$ws.Range("H77").Value = 31251492
$ws.Range("I77").Value = 35715260
$ws.Range("J77").Value = 5124.25
$ws.Range("K77").Value = 178576300
$ws.Range("L77").Value = 25621.25
$ws.Range("M77").Value = -178571932
$ws.Range("N77").Value = -34357.25

$ws.Range("H122").Value = 19611888
$ws.Range("I122").Value = 4290.5386
$ws.Range("J122").Value = 83336580
$ws.Range("K122").Value = 12871.6158
$ws.Range("L122").Value = 250009740
$ws.Range("M122").Value = -10421.6158
$ws.Range("N122").Value = -250014640

$ws.Range("H132").Value = 34577224
$ws.Range("I132").Value = 17056
$ws.Range("J132").Value = 91130220
$ws.Range("K132").Value = 51168
$ws.Range("L132").Value = 273390660
$ws.Range("M132").Value = -48638
$ws.Range("N132").Value = -273395720

$ws.Range("H137").Value = 50321.332
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 50321.332
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 50321.332
$ws.Range("N137").Value = -60521.332

$ws.Range("H139").Value = 174999
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 174999
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 174999
$ws.Range("N139").Value = -185279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 48140
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 48140
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 48140
$ws.Range("N58").Value = -48728

$ws.Range("H59").Value = 71323
$ws.Range("I59").Value = 51985

$ws.Range("H134").Value = 3511.25
$ws.Range("I134").Value = 3701.389
$ws.Range("J134").Value = 1800
$ws.Range("K134").Value = 11104.167
$ws.Range("L134").Value = 5400
$ws.Range("M134").Value = -8569.167000000001
$ws.Range("N134").Value = -10470

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3467.25
$ws.Range("I16").Value = 3138.5
$ws.Range("J16").Value = 4124.75
$ws.Range("K16").Value = 3138.5
$ws.Range("L16").Value = 4124.75
$ws.Range("M16").Value = -2851.5
$ws.Range("N16").Value = -4698.75

$ws.Range("H31").Value = 10004546
$ws.Range("I31").Value = 3029.625
$ws.Range("J31").Value = 27785020
$ws.Range("K31").Value = 3029.625
$ws.Range("L31").Value = 27785020
$ws.Range("M31").Value = -2734.625
$ws.Range("N31").Value = -27785610

$ws.Range("H34").Value = 10004546
$ws.Range("I34").Value = 3029.625
$ws.Range("J34").Value = 27785020
$ws.Range("K34").Value = 3029.625
$ws.Range("L34").Value = 27785020
$ws.Range("M34").Value = -2827.625
$ws.Range("N34").Value = -27785424

$ws.Range("H86").Value = 2390.4443
$ws.Range("I86").Value = 2144.8572
$ws.Range("J86").Value = 3250
$ws.Range("K86").Value = 2144.8572
$ws.Range("L86").Value = 3250
$ws.Range("M86").Value = -1021.8572
$ws.Range("N86").Value = -5496

$ws.Range("H89").Value = 2390.4443
$ws.Range("I89").Value = 2144.8572
$ws.Range("J89").Value = 3250
$ws.Range("K89").Value = 10724.286
$ws.Range("L89").Value = 16250
$ws.Range("M89").Value = -5108.286
$ws.Range("N89").Value = -27482

$ws.Range("H113").Value = 3467.25
$ws.Range("I113").Value = 3138.5
$ws.Range("J113").Value = 4124.75
$ws.Range("K113").Value = 3138.5
$ws.Range("L113").Value = 4124.75
$ws.Range("M113").Value = -968.5
$ws.Range("N113").Value = -8464.75

$ws.Range("H132").Value = 170833.75
$ws.Range("I132").Value = 670004
$ws.Range("J132").Value = 4443.6665
$ws.Range("K132").Value = 2010012
$ws.Range("L132").Value = 13330.9995
$ws.Range("M132").Value = -2007482
$ws.Range("N132").Value = -18390.9995

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H141").Value = 127498.164
$ws.Range("I141").Value = 30000
$ws.Range("J141").Value = 146997.8
$ws.Range("K141").Value = 30000
$ws.Range("L141").Value = 146997.8
$ws.Range("M141").Value = -24820
$ws.Range("N141").Value = -157357.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1492.6666
$ws.Range("I5").Value = 343.125
$ws.Range("J5").Value = 2806.4285
$ws.Range("K5").Value = 1029.375
$ws.Range("L5").Value = 8419.2855
$ws.Range("M5").Value = -917.375
$ws.Range("N5").Value = -8643.2855

$ws.Range("H50").Value = 1250.6666
$ws.Range("I50").Value = 1168.3334
$ws.Range("J50").Value = 1333
$ws.Range("K50").Value = 3505.0002
$ws.Range("L50").Value = 3999
$ws.Range("M50").Value = -3024.0002
$ws.Range("N50").Value = -4961

$ws.Range("H53").Value = 1250.6666
$ws.Range("I53").Value = 1168.3334
$ws.Range("J53").Value = 1333
$ws.Range("K53").Value = 3505.0002
$ws.Range("L53").Value = 3999
$ws.Range("M53").Value = -3024.0002
$ws.Range("N53").Value = -4961

$ws.Range("H128").Value = 200000
$ws.Range("I128").Value = 200000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 600000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -595020

$ws.Range("H132").Value = 4448781.5
$ws.Range("I132").Value = 2121
$ws.Range("J132").Value = 11118772
$ws.Range("K132").Value = 19089
$ws.Range("L132").Value = 100068948
$ws.Range("M132").Value = -16559
$ws.Range("N132").Value = -100074008

$ws.Range("H135").Value = 1492.6666
$ws.Range("I135").Value = 343.125
$ws.Range("J135").Value = 2806.4285
$ws.Range("K135").Value = 3088.125
$ws.Range("L135").Value = 25257.8565
$ws.Range("M135").Value = -553.125
$ws.Range("N135").Value = -30327.8565

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 45000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 45000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 45000
$ws.Range("N33").Value = -45504
$ws.Range("M33").ClearContents()

$ws.Range("H102").Value = 1168.1666
$ws.Range("I102").Value = 1002
$ws.Range("J102").Value = 1999
$ws.Range("K102").Value = 1002
$ws.Range("L102").Value = 1999
$ws.Range("M102").Value = 620
$ws.Range("N102").Value = -5243

$ws.Range("H132").Value = 4375.875
$ws.Range("I132").Value = 4668
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 14004
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -11474
$ws.Range("N132").Value = -15558.5

$ws.Range("H133").Value = 197998
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 197998
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 197998
$ws.Range("N133").Value = -208118

$ws.Range("H135").Value = 119988.6
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 119988.6
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 119988.6
$ws.Range("N135").Value = -130128.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1130.5
$ws.Range("I22").Value = 908.4666999999999
$ws.Range("J22").Value = 1386.6923
$ws.Range("K22").Value = 908.4666999999999
$ws.Range("L22").Value = 1386.6923
$ws.Range("M22").Value = -613.4666999999999
$ws.Range("N22").Value = -1976.6923

$ws.Range("H27").Value = 1130.5
$ws.Range("I27").Value = 908.4666999999999
$ws.Range("J27").Value = 1386.6923
$ws.Range("K27").Value = 908.4666999999999
$ws.Range("L27").Value = 1386.6923
$ws.Range("M27").Value = -801.4666999999999
$ws.Range("N27").Value = -1600.6923

$ws.Range("H40").Value = 5283.04
$ws.Range("I40").Value = 5158.5454
$ws.Range("J40").Value = 6196
$ws.Range("K40").Value = 5158.5454
$ws.Range("L40").Value = 6196
$ws.Range("M40").Value = -5022.5454
$ws.Range("N40").Value = -6468

$ws.Range("H104").Value = 21934.75
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 21934.75
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 21934.75
$ws.Range("N104").Value = -28922.75

$ws.Range("H122").Value = 3382388
$ws.Range("I122").Value = 3935.2
$ws.Range("J122").Value = 10420831
$ws.Range("K122").Value = 11805.6
$ws.Range("L122").Value = 31262493
$ws.Range("M122").Value = -9355.599999999999
$ws.Range("N122").Value = -31267393

$ws.Range("H131").Value = 76853.5
$ws.Range("I131").Value = 52648
$ws.Range("J131").Value = 88956.25
$ws.Range("K131").Value = 52648
$ws.Range("L131").Value = 88956.25
$ws.Range("M131").Value = -47608
$ws.Range("N131").Value = -99036.25

$ws.Range("H140").Value = 52291.5
$ws.Range("I140").Value = 40390
$ws.Range("J140").Value = 87996
$ws.Range("K140").Value = 40390
$ws.Range("L140").Value = 87996
$ws.Range("M140").Value = -35210
$ws.Range("N140").Value = -98356

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 375050
$ws.Range("I2").Value = 750000

$ws.Range("H41").Value = 1000000000
$ws.Range("I41").Value = 1000000000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1000000000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -999999610
$ws.Range("N41").ClearContents()

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H122").Value = 10529706
$ws.Range("I122").Value = 1795.2
$ws.Range("J122").Value = 22227384
$ws.Range("K122").Value = 5385.6
$ws.Range("L122").Value = 66682152
$ws.Range("M122").Value = -2935.6
$ws.Range("N122").Value = -66687052

$ws.Range("H126").Value = 10428088
$ws.Range("I126").Value = 12831523
$ws.Range("J126").Value = 13200
$ws.Range("K126").Value = 38494569
$ws.Range("L126").Value = 39600
$ws.Range("M126").Value = -38492099
$ws.Range("N126").Value = -44540

$ws.Range("H141").Value = 68497.664
$ws.Range("I141").Value = 109000
$ws.Range("J141").Value = 63434.875
$ws.Range("K141").Value = 109000
$ws.Range("L141").Value = 63434.875
$ws.Range("M141").Value = -94255.39999999999
$ws.Range("N141").Value = -73794.875
